$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap data between row 41 and row 42 (columns F:V) ---
$row41 = $ws.Range("F41:V41").Value()
$row42 = $ws.Range("F42:V42").Value()
$ws.Range("F41:V41").Value = $row42
$ws.Range("F42:V42").Value = $row41

# --- Append 3 new match rows (82, 83, 84) ---
# Copy formatting (styles) from the last existing row (81) down to the new rows
$ws.Range("A81:V81").Copy()
$ws.Range("A82:V84").PasteSpecial(-4122)

# Row 82
$ws.Range("A82").Value = 81
$ws.Range("B82").Value = "poland"
$ws.Range("C82").Value = "division-2"
$ws.Range("D82").Value = "2023-2024"
$ws.Range("E82").Value = 45191.70833333334
$ws.Range("F82").Value = "Wisla Pulawy"
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = "Lech Poznan II"
$ws.Range("I82").Value = 2
$ws.Range("J82").Value = 1.61
$ws.Range("K82").Value = "21/09/2023 04:12"
$ws.Range("L82").Value = 1.68
$ws.Range("M82").Value = "22/09/2023 16:50"
$ws.Range("N82").Value = 3.74
$ws.Range("O82").Value = "21/09/2023 04:12"
$ws.Range("P82").Value = 4.02
$ws.Range("Q82").Value = "22/09/2023 16:50"
$ws.Range("R82").Value = 4.37
$ws.Range("S82").Value = "21/09/2023 04:12"
$ws.Range("T82").Value = 4.41
$ws.Range("U82").Value = "22/09/2023 16:50"
$ws.Range("V82").Value = "https://www.betexplorer.com/football/poland/division-2/wisla-pulawy-lech-poznan/QuRUSzJq/"

# Row 83
$ws.Range("A83").Value = 82
$ws.Range("B83").Value = "poland"
$ws.Range("C83").Value = "division-2"
$ws.Range("D83").Value = "2023-2024"
$ws.Range("E83").Value = 45191.8125
$ws.Range("F83").Value = "Chojniczanka"
$ws.Range("G83").Value = 5
$ws.Range("H83").Value = "Stezyca"
$ws.Range("I83").Value = 1
$ws.Range("J83").Value = 2.08
$ws.Range("K83").Value = "21/09/2023 06:42"
$ws.Range("L83").Value = 2.12
$ws.Range("M83").Value = "22/09/2023 19:20"
$ws.Range("N83").Value = 3.2
$ws.Range("O83").Value = "21/09/2023 06:42"
$ws.Range("P83").Value = 3.31
$ws.Range("Q83").Value = "22/09/2023 17:31"
$ws.Range("R83").Value = 3.13
$ws.Range("S83").Value = "21/09/2023 06:42"
$ws.Range("T83").Value = 3.36
$ws.Range("U83").Value = "22/09/2023 19:20"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/poland/division-2/chojniczanka-stezyca/YsnjzJQ9/"

# Row 84
$ws.Range("A84").Value = 83
$ws.Range("B84").Value = "poland"
$ws.Range("C84").Value = "division-2"
$ws.Range("D84").Value = "2023-2024"
$ws.Range("E84").Value = 45191.83333333334
$ws.Range("F84").Value = "KKS Kalisz"
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = "Hutnik Krakow"
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1.95
$ws.Range("K84").Value = "21/09/2023 07:12"
$ws.Range("L84").Value = 2.02
$ws.Range("M84").Value = "22/09/2023 19:39"
$ws.Range("N84").Value = 3.28
$ws.Range("O84").Value = "21/09/2023 07:12"
$ws.Range("P84").Value = 3.38
$ws.Range("Q84").Value = "22/09/2023 19:39"
$ws.Range("R84").Value = 3.34
$ws.Range("S84").Value = "21/09/2023 07:12"
$ws.Range("T84").Value = 3.56
$ws.Range("U84").Value = "22/09/2023 19:39"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/poland/division-2/kks-kalisz-hutnik-krakow/CY6PlIu3/"
